# Sync attendance_reports: swap the order of names in the "Recorded By"
# column (G) from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row where that exact value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = 7
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
    }
}
